$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two existing hyperlinks (on E9 and E10) need to end up one row lower
# (E10 and E11) once the new row is inserted above them. This host does not
# shift hyperlink anchors automatically on row insert, and Hyperlink.Range is
# not writable, so remove them first and recreate them afterwards.
$ws.Range("E9").Hyperlinks.Delete()
$ws.Range("E10").Hyperlinks.Delete()

# Insert a new row above row 5 ("qui_es_tu"), pushing everything from the old
# row 5 onward down by one.
$ws.Range("A5").EntireRow.Insert()

# Fill in the new row 5 with the "phrase_feedback" entry.
$ws.Range("A5").Value = "phrase_feedback"
$ws.Range("D5").Value = "Texte"
$ws.Range("E5").Value = "Bonjour, cet ID correspond au message envoyé quand un utilisateur donne son avis sur une réponse."
$ws.Range("G5").Value = "Phrase feedback.; J'ai laissé mon avis !"

# The insert pushes the row-below formatting onto the new blank cells (G5
# inherited the shaded "réponse" style, E5 correctly kept the wrap style).
# Reset G5 back to the default/unstyled look used for a plain text note.
$ws.Range("G5").Style = "Normal"

# This row's text wraps onto two lines like the other multi-line rows.
$ws.Rows("5:5").RowHeight = 28.5

# Recreate the hyperlinks one row further down than before.
$ws.Hyperlinks.Add($ws.Range("E10"), "https://i.imgur.com/nGF1K8f.jpg")
$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.google.fr")

# Adding a hyperlink re-styles the cell with the default "followed link"
# look; restore the sheet's own hyperlink style (wrapped text) instead.
$ws.Range("E10").Style = "Lien hypertexte"
$ws.Range("E10").WrapText = $true
$ws.Range("E11").Style = "Lien hypertexte"
$ws.Range("E11").WrapText = $true

# Match the author's final selection.
$ws.Range("E5").Select()
